$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$nm.Delete()
$nm2 = $p.NotesMaster
Write-Output "done"
